# Being Creative Research Project - stimuli_info_sheet.xlsx edit
#
# Adds two new "ai image info" prompts (for the bottle_and_glass and
# geometric_dog stimuli rows), grows the affected rows to fit the newly
# wrapped text, and updates the sheet's selection/scroll position to
# where the user ended up after editing (around row 9/10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stimuli List")

# --- New content: D7 (bottle_and_glass) and D9 (geometric_dog) "ai image info" ---
$ws.Range("D7").Value = "A minimalist line illustration of a bottle with a glass nested inside, creating a simple yet clever composition. This artwork captures the beauty of modern minimalism, perfect for those who appreciate refined, abstract line art.`nhow creative and abstract can you be? focus on it being inside"

$ws.Range("D9").Value = "This minimalist line art illustration features a geometric interpretation of a dog, blending simplicity and precision. Created with clean lines and abstract shapes, it highlights the beauty of symmetry and modern design. Perfect for pet lovers, contemporary art collections, and decor enthusiasts, this artwork embodies elegance and creativity while celebrating canine charm. "

# --- Row heights grow to accommodate newly wrapped / longer text ---
$ws.Rows.Item(4).RowHeight = 131.15
$ws.Rows.Item(5).RowHeight = 262.3
$ws.Rows.Item(6).RowHeight = 131.15
$ws.Rows.Item(7).RowHeight = 204
$ws.Rows.Item(8).RowHeight = 174.9
$ws.Rows.Item(9).RowHeight = 247.75

# --- Update view: scrolled down with C10 selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("C10").Select()
